$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8 corresponds to year 2025 in metricas_recorrencia_anual.
# Update total_customers, returning_customers, new_customers and the
# derived rate columns to the refreshed figures from the commit.
$ws.Range("C8").Value = 1009
$ws.Range("D8").Value = 165
$ws.Range("E8").Value = 844
$ws.Range("F8").Value = 6.767842493847416
$ws.Range("G8").Value = 83.64717542120911
$ws.Range("H8").Value = 16.35282457879088
